# "Ajout de methodes et update backings"
#
# Adds a handful of new Q/A rows to Feuil1 (sheet1) documenting a few more
# DAO/manager methods, and nudges the saved viewport/selection on Feuil1 and
# Feuil2 to where the author last left them.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- New content on Feuil1 -------------------------------------------------
# New shared strings must be introduced in the same order the author typed
# them so they land on the expected shared-string indices:
#   140 get membre by nickname(nickname)
#   141 get membre by nickname et password(nickname+password)
#   142 MessageManager.deleteMessage(Message message);
#   143 getClinsdoeilRecus(int TOId){
#   144 arraylist clinsdoeils qui un user a recu!
#   145 ClinsdoeilManager.deleteClindoeil(Clinsdoeil clin){

$ws1.Range("A69").Value = "get membre by nickname(nickname)"
$ws1.Range("A69").Font.Color = 255

$ws1.Range("A70").Value = "get membre by nickname et password(nickname+password)"
$ws1.Range("A70").Font.Color = 255

$ws1.Range("A68").Value = "MessageManager.deleteMessage(Message message);"
$ws1.Range("A68").Font.Color = 255

$ws1.Range("A71").Value = "getClinsdoeilRecus(int TOId){"
$ws1.Range("A71").Font.Color = 255

$ws1.Range("B71").Value = "arraylist clinsdoeils qui un user a recu!"
$ws1.Range("B71").Font.Color = 255

$ws1.Range("A72").Value = "ClinsdoeilManager.deleteClindoeil(Clinsdoeil clin){"

# --- Update the saved selection / viewport ---------------------------------
# Feuil2 first so that re-activating Feuil1 afterwards leaves it as the
# selected tab, matching the original workbook.
$ws2.Activate()
$ws2.Range("B30").Select()

$ws1.Activate()
$ws1.Range("A72").Select()
